$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "64.247.66"
$ws.Range("E2").Value = "  +1.76%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.627.08"
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "602.80"
$ws.Range("E5").Value = "  +0.27%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "150.93"
$ws.Range("E6").Value = "  +3.28%  "
$ws.Range("E7").Value = "  -0.02%  "
$ws.Range("E9").Value = "  +2.39%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "5.72"
$ws.Range("E10").Value = "  +2.18%  "
$ws.Range("E11").Value = "  +6.66%  "
$ws.Range("E12").Value = "  -0.68%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "27.73"
$ws.Range("E13").Value = "  +2.25%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "3.099.74"
$ws.Range("E14").Value = "  +0.01%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "64.092.27"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000150"
$ws.Range("E16").Value = "  +3.95%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.648.35"
$ws.Range("E17").Value = "  +0.70%  "
$ws.Range("E18").Value = "  +8.10%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.67"
$ws.Range("E19").Value = "  +4.16%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "351.88"
$ws.Range("E20").Value = "  +3.55%  "
$ws.Range("E21").Value = "  +1.68%  "
$ws.Range("E22").Value = "  +0.03%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.73"
$ws.Range("E23").Value = "  +2.86%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "66.65"
$ws.Range("E24").Value = "  +0.28%  "
$ws.Range("E25").Value = "  +17.02%  "
$ws.Range("E26").Value = "  +4.48%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.26"
$ws.Range("E27").Value = "  +7.24%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.166"
$ws.Range("E28").Value = "  +1.87%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.18"
$ws.Range("E29").Value = "  +4.41%  "
$ws.Range("E30").Value = "  -0.12%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "539.24"
$ws.Range("E31").Value = "  -1.66%  "
$ws.Range("E32").Value = "  +3.50%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0₃0858"
$ws.Range("E33").Value = "  +7.10%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.76"
$ws.Range("E34").Value = "  +1.02%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.31"
$ws.Range("E35").Value = "  +0.34%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "167.91"
$ws.Range("E36").Value = "  +1.26%  "
$ws.Range("E37").Value = "  +7.92%  "
$ws.Range("E38").Value = "  +1.93%  "
$ws.Range("E39").Value = "  -0.09%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "19.52"
$ws.Range("E40").Value = "  +3.05%  "
$ws.Range("B41").Value = "Aave"
$ws.Range("C41").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "170.39"
$ws.Range("E41").Value = "  +1.74%  "
$ws.Range("B42").Value = "USDe"
$ws.Range("C42").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.999"
$ws.Range("E42").Value = "  +0.03%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "39.98"
$ws.Range("E43").Value = "  +0.97%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.95"
$ws.Range("E44").Value = "  +6.01%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0594"
$ws.Range("E45").Value = "  +4.84%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "21.61"
$ws.Range("E46").Value = "  -3.03%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.631"
$ws.Range("E47").Value = "  +1.48%  "
$ws.Range("E48").Value = "  +2.03%  "
$ws.Range("E49").Value = "  +14.11%  "
$ws.Range("E50").Value = "  +1.14%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "19.34"
$ws.Range("E51").Value = "  +4.21%  "

Write-Host "Applied cryptos update"
